# Weekly Status Report fixes ("fixing documentation so it matches")
#
# This script rewrites the contents of seven table cells so that the
# document matches the target revision. Word's Find/Replace in this
# runtime searches the whole story regardless of the Range it is called
# on, so every edit below is instead done with Range.InsertXML, which
# DOES respect the Range it is invoked on (but replaces the whole
# paragraph that Range lives in) -- so each payload reconstructs the
# complete paragraph (every run, changed or not) for the target cell.

$d = $word.ActiveDocument

function Set-CellXml($cell, [string]$innerXml) {
    $rng = $d.Range($cell.Range.Start, $cell.Range.End)
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p>' + $innerXml + '</w:p></w:body>' +
           '</w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($pkg)
}

$RPR = '<w:rPr><w:rFonts w:eastAsia="Calibri"/><w:color w:val="000000"/></w:rPr>'

# ---------------------------------------------------------------------
# Hunk 1: Table 2 ("Tasks Completed This Period"), row 3 (Harman),
# Status cell: the old "_GoBack" bookmark used to split "Co" / "mpleted "
# into two runs; the new text is a single clean run again.
# ---------------------------------------------------------------------
$t2 = $d.Tables.Item(2)
$cellHarmanStatus = $t2.Rows.Item(3).Cells.Item(5)
$xml1 = '<w:pPr><w:rPr><w:rFonts w:eastAsia="Calibri"/><w:color w:val="000000"/></w:rPr></w:pPr>' +
        '<w:r>' + $RPR + '<w:t xml:space="preserve">Completed </w:t></w:r>'
Set-CellXml $cellHarmanStatus $xml1

# ---------------------------------------------------------------------
# Hunk 2: Table 2, row 5 (John), Status cell: "Completed " -> "In progress "
# (split across two runs to match the target).
# ---------------------------------------------------------------------
$cellJohnStatus = $t2.Rows.Item(5).Cells.Item(5)
$xml2 = '<w:pPr><w:rPr><w:rFonts w:eastAsia="Calibri"/><w:color w:val="000000"/></w:rPr></w:pPr>' +
        '<w:r>' + $RPR + '<w:t>In progress</w:t></w:r>' +
        '<w:r>' + $RPR + '<w:t xml:space="preserve"> </w:t></w:r>'
Set-CellXml $cellJohnStatus $xml2

# ---------------------------------------------------------------------
# Hunk 3: Table 3 ("Tasks Planned but Not Completed" / in progress),
# row 2 (Henry), Task Description cell: "Backend for Admin setup" ->
# "Backend for admin setup", with the last-edit "_GoBack" bookmark now
# sitting at this (new) split point between "a" and "dmin setup".
# ---------------------------------------------------------------------
$t3 = $d.Tables.Item(3)
$cellHenryInProgress = $t3.Rows.Item(2).Cells.Item(2)
$xml3 = '<w:pPr><w:rPr><w:rFonts w:eastAsia="Calibri"/><w:i/><w:color w:val="5B9BD5"/></w:rPr></w:pPr>' +
        '<w:r>' + $RPR + '<w:t xml:space="preserve">Backend for </w:t></w:r>' +
        '<w:r>' + $RPR + '<w:t>a</w:t></w:r>' +
        '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
        '<w:r>' + $RPR + '<w:t>dmin setup</w:t></w:r>'
Set-CellXml $cellHenryInProgress $xml3

# ---------------------------------------------------------------------
# Hunk 4: Table 4 ("Tasks Planned for Next Period"), row 2 (Henry),
# Task Description cell.
# ---------------------------------------------------------------------
$t4 = $d.Tables.Item(4)
$cellHenryPlanned = $t4.Rows.Item(2).Cells.Item(2)
$xml4 = '<w:pPr><w:rPr><w:rFonts w:eastAsia="Calibri"/><w:color w:val="000000"/></w:rPr></w:pPr>' +
        '<w:r>' + $RPR + '<w:t>Finish in progress tasks,</w:t></w:r>' +
        '<w:r>' + $RPR + '<w:t xml:space="preserve"> work on any backend changes related to a</w:t></w:r>' +
        '<w:r>' + $RPR + '<w:t>dmin setup, customer documentation + testing</w:t></w:r>' +
        '<w:r>' + $RPR + '<w:t xml:space="preserve">  </w:t></w:r>'
Set-CellXml $cellHenryPlanned $xml4

# ---------------------------------------------------------------------
# Hunk 5: Table 4, row 3 (Harman), Task Description cell.
# ---------------------------------------------------------------------
$cellHarmanPlanned = $t4.Rows.Item(3).Cells.Item(2)
$xml5 = '<w:pPr><w:rPr><w:rFonts w:eastAsia="Calibri"/><w:color w:val="000000"/></w:rPr></w:pPr>' +
        '<w:r>' + $RPR + '<w:t>Finish in progress tasks, set up receiving/sending KPI for a</w:t></w:r>' +
        '<w:r>' + $RPR + '<w:t xml:space="preserve">dmin, customer documentation + testing  </w:t></w:r>'
Set-CellXml $cellHarmanPlanned $xml5

# ---------------------------------------------------------------------
# Hunk 6: Table 4, row 4 (John), Task Description cell. The trailing
# "GUI cleanup" + " + testing  " runs are untouched by the diff, so they
# are carried forward unchanged.
# ---------------------------------------------------------------------
$cellJohnPlanned = $t4.Rows.Item(4).Cells.Item(2)
$xml6 = '<w:pPr><w:rPr><w:rFonts w:eastAsia="Calibri"/><w:color w:val="000000"/></w:rPr></w:pPr>' +
        '<w:r>' + $RPR + '<w:t>Finish in progress tasks,</w:t></w:r>' +
        '<w:r>' + $RPR + '<w:t xml:space="preserve"> help with a</w:t></w:r>' +
        '<w:r>' + $RPR + '<w:t xml:space="preserve">dmin setup, </w:t></w:r>' +
        '<w:r>' + $RPR + '<w:t>GUI cleanup</w:t></w:r>' +
        '<w:r>' + $RPR + '<w:t xml:space="preserve"> + testing  </w:t></w:r>'
Set-CellXml $cellJohnPlanned $xml6

# ---------------------------------------------------------------------
# Hunk 7: Table 4, row 5 (Oliver), Task Description cell.
# ---------------------------------------------------------------------
$cellOliverPlanned = $t4.Rows.Item(5).Cells.Item(2)
$xml7 = '<w:pPr><w:rPr><w:rFonts w:eastAsia="Calibri"/><w:color w:val="000000"/></w:rPr></w:pPr>' +
        '<w:r>' + $RPR + '<w:t>Finish in progress tasks,</w:t></w:r>' +
        '<w:r>' + $RPR + '<w:t xml:space="preserve"> help with a</w:t></w:r>' +
        '<w:r>' + $RPR + '<w:t xml:space="preserve">dmin setup, customer documentation + testing  </w:t></w:r>'
Set-CellXml $cellOliverPlanned $xml7

Write-Output "done"
